# Update the "想去人数" (want-to-go count) figures on the sheets that hold
# the event data: "展览" and "全部类型" (sheet1 and sheet4 both carry the
# same table). "演出" and "本地生活" only contain header rows, so they are
# left untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 65
    $ws.Range("F3").Value = 1300
    $ws.Range("F5").Value = 5
    $ws.Range("F8").Value = 171
}
